$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()

$find.Execute(
    "para o dia 20, às 11 horas, na EST, Lab Internet of things.",
    $true,
    $false,
    $false,
    $false,
    $false,
    $false,
    1,
    $false,
    "para o dia 20 de Setembro, às 11 horas, na EST, Lab Internet of things.",
    2
)
